$d = $word.ActiveDocument

# Literal separator line used as a thin visual divider between sections.
$sepLine = "────────────────────────────────────────────────────────────"

# Collect every paragraph that must be removed:
#   1) the thin "before=40" spacer paragraphs that sit right after a table
#   2) the "────..." separator paragraphs (color CCCCCC, sz 16)
$targets = @()
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq $sepLine) {
        $targets += $p
    } elseif ($txt -eq "" -and $p.Range.SpaceBefore -eq 2 -and $p.Range.InlineShapes.Count -eq 0) {
        $targets += $p
    }
}

# Delete from the end of the document backwards so earlier ranges in the
# list are never invalidated/shifted by a later delete.
for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $targets[$i].Range.Delete()
}
